{"js": "// Word Homework Chapter 01 edit:\n//  1) 1.B prompt: re-type the straight quotes around \"goed\" / \"went\" so each\n//     quote mark and quoted word lands in its own run (no visible text change).\n//  2) The three dialogue lines (A: / B: / A:) lose the straight quotes that\n//     wrapped the quoted speech, splitting the remainder into a leading-space\n//     run and a quote-free text run.\n//\n// Helper: wrap a raw <w:p>...</w:p> fragment in the minimal OOXML package\n// envelope that Range.insertOoxml() expects.\nfunction wrapParagraphOoxml(pInnerXml) {\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' + pInnerXml + '</w:p></w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\n// Find the unique range in the paragraph whose text equals `oldText`, then\n// replace just that range with the run-split OOXML given in `runsXml`\n// (a concatenation of <w:r>...</w:r> elements). Leaves any other runs in\n// the paragraph (e.g. the bold \"1.B  \" / italic \"A:\" label run) untouched.\nasync function splitRun(context, paragraph, oldText, runsXml) {\n  const results = paragraph.search(oldText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"splitRun: text not found: \" + oldText);\n  }\n\n  const target = results.items[0];\n  target.insertOoxml(wrapParagraphOoxml(runsXml), \"Replace\");\n  await context.sync();\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet para1B = null;\nlet paraA1 = null; // \"A: <q>You going to the thing tonight?<q>\"\nlet paraB = null;  // \"B: <q>Bro, I'm exhausted.<q>\"\nlet paraA2 = null; // \"A: <q>Bet.<q>\"\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (para1B === null && t.indexOf('1.B') === 0 && t.indexOf('\"goed\"') !== -1) {\n    para1B = paragraphs.items[i];\n  } else if (paraA1 === null && t === 'A: \"You going to the thing tonight?\"') {\n    paraA1 = paragraphs.items[i];\n  } else if (paraB === null && t === 'B: \"Bro, I\\'m exhausted.\"') {\n    paraB = paragraphs.items[i];\n  } else if (paraA2 === null && t === 'A: \"Bet.\"') {\n    paraA2 = paragraphs.items[i];\n  }\n}\n\nif (!para1B || !paraA1 || !paraB || !paraA2) {\n  throw new Error(\n    \"Could not locate all target paragraphs (1B=\" + !!para1B +\n    \" A1=\" + !!paraA1 + \" B=\" + !!paraB + \" A2=\" + !!paraA2 + \")\"\n  );\n}\n\n// 1) 1.B \u2014 split the single answer run into 9 runs around the re-typed quotes.\nawait splitRun(\n  context,\n  para1B,\n  ' The chapter argues that children who say \"goed\" instead of \"went\" are actually demonstrating sophisticated language learning. Explain why linguists see this as evidence of productivity rather than error. (1 paragraph)',\n  '<w:r><w:t xml:space=\"preserve\"> The chapter argues that children who say </w:t></w:r>' +\n  '<w:r><w:t>&quot;</w:t></w:r>' +\n  '<w:r><w:t>goed</w:t></w:r>' +\n  '<w:r><w:t>&quot;</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> instead of </w:t></w:r>' +\n  '<w:r><w:t>&quot;</w:t></w:r>' +\n  '<w:r><w:t>went</w:t></w:r>' +\n  '<w:r><w:t>&quot;</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> are actually demonstrating sophisticated language learning. Explain why linguists see this as evidence of productivity rather than error. (1 paragraph)</w:t></w:r>'\n);\n\n// 2) A: \"You going to the thing tonight?\" -> A: You going to the thing tonight?\nawait splitRun(\n  context,\n  paraA1,\n  ' \"You going to the thing tonight?\"',\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t>You going to the thing tonight?</w:t></w:r>'\n);\n\n// 3) B: \"Bro, I'm exhausted.\" -> B: Bro, I'm exhausted.\nawait splitRun(\n  context,\n  paraB,\n  ' \"Bro, I\\'m exhausted.\"',\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t>Bro, I&apos;m exhausted.</w:t></w:r>'\n);\n\n// 4) A: \"Bet.\" -> A: Bet.\nawait splitRun(\n  context,\n  paraA2,\n  ' \"Bet.\"',\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t>Bet.</w:t></w:r>'\n);\n", "ps1": "# Word Homework Chapter 01 edit:\n#  1) 1.B prompt: re-type the straight quotes around \"goed\" / \"went\" so each\n#     quote mark and quoted word lands in its own run (no visible text change).\n#  2) The three dialogue lines (A: / B: / A:) lose the straight quotes that\n#     wrapped the quoted speech, splitting the remainder into a leading-space\n#     run and a quote-free text run.\n\nfunction Wrap-ParagraphOoxml {\n    param([string]$PInnerXml)\n    return '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n           '<pkg:part pkg:name=\"/word/document.xml\" ' +\n           'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n           '<pkg:xmlData>' +\n           '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n           '<w:body><w:p>' + $PInnerXml + '</w:p></w:body>' +\n           '</w:document>' +\n           '</pkg:xmlData></pkg:part></pkg:package>'\n}\n\n# Pull this paragraph's own <w:pPr>...</w:pPr> (if any) out of its WordOpenXML\n# so the rebuilt <w:p> keeps its original indentation/spacing.\nfunction Get-ParagraphPPr {\n    param($Paragraph)\n    $full = $Paragraph.Range.WordOpenXML\n    if ($full -match '<w:p\\b[^>]*>(.*?)</w:p>') {\n        $inner = $matches[1]\n        if ($inner -match '<w:pPr\\b[^>]*>.*?</w:pPr>|<w:pPr\\b[^>]*/>') {\n            return $matches[0]\n        }\n    }\n    return \"\"\n}\n\n# Replace this paragraph's whole content (all its runs) with a fresh run-split\n# described by $RunsXml (a concatenation of <w:r>...</w:r> elements), keeping\n# the paragraph's own formatting (<w:pPr>) intact.\nfunction Set-ParagraphRuns {\n    param($Paragraph, [string]$RunsXml)\n    $pPr = Get-ParagraphPPr $Paragraph\n    $ooxml = Wrap-ParagraphOoxml ($pPr + $RunsXml)\n    $Paragraph.Range.InsertXML($ooxml) | Out-Null\n}\n\n$d = $word.ActiveDocument\n\n$para1B = $null\n$paraA1 = $null   # A: \"You going to the thing tonight?\"\n$paraB  = $null   # B: \"Bro, I'm exhausted.\"\n$paraA2 = $null   # A: \"Bet.\"\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if (($null -eq $para1B) -and $t.StartsWith('1.B') -and ($t.IndexOf('\"goed\"') -ge 0)) {\n        $para1B = $p\n    } elseif (($null -eq $paraA1) -and ($t -eq 'A: \"You going to the thing tonight?\"' + [char]13)) {\n        $paraA1 = $p\n    } elseif (($null -eq $paraB) -and ($t -eq \"B: `\"Bro, I'm exhausted.`\"\" + [char]13)) {\n        $paraB = $p\n    } elseif (($null -eq $paraA2) -and ($t -eq 'A: \"Bet.\"' + [char]13)) {\n        $paraA2 = $p\n    }\n}\n\nif (($null -eq $para1B) -or ($null -eq $paraA1) -or ($null -eq $paraB) -or ($null -eq $paraA2)) {\n    throw \"Could not locate all target paragraphs\"\n}\n\n# 1) 1.B \u2014 split the single answer run into 9 runs around the re-typed quotes.\n$runs1B = '<w:r><w:t xml:space=\"preserve\"> The chapter argues that children who say </w:t></w:r>' +\n          '<w:r><w:t>&quot;</w:t></w:r>' +\n          '<w:r><w:t>goed</w:t></w:r>' +\n          '<w:r><w:t>&quot;</w:t></w:r>' +\n          '<w:r><w:t xml:space=\"preserve\"> instead of </w:t></w:r>' +\n          '<w:r><w:t>&quot;</w:t></w:r>' +\n          '<w:r><w:t>went</w:t></w:r>' +\n          '<w:r><w:t>&quot;</w:t></w:r>' +\n          '<w:r><w:t xml:space=\"preserve\"> are actually demonstrating sophisticated language learning. Explain why linguists see this as evidence of productivity rather than error. (1 paragraph)</w:t></w:r>'\n\n$bold1B = '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">1.B  </w:t></w:r>'\nSet-ParagraphRuns $para1B ($bold1B + $runs1B)\n\n# 2) A: \"You going to the thing tonight?\" -> A: You going to the thing tonight?\n$runsA1 = '<w:r><w:rPr><w:i/></w:rPr><w:t>A:</w:t></w:r>' +\n          '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n          '<w:r><w:t>You going to the thing tonight?</w:t></w:r>'\nSet-ParagraphRuns $paraA1 $runsA1\n\n# 3) B: \"Bro, I'm exhausted.\" -> B: Bro, I'm exhausted.\n$runsB = '<w:r><w:rPr><w:i/></w:rPr><w:t>B:</w:t></w:r>' +\n         '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n         '<w:r><w:t>Bro, I&apos;m exhausted.</w:t></w:r>'\nSet-ParagraphRuns $paraB $runsB\n\n# 4) A: \"Bet.\" -> A: Bet.\n$runsA2 = '<w:r><w:rPr><w:i/></w:rPr><w:t>A:</w:t></w:r>' +\n          '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n          '<w:r><w:t>Bet.</w:t></w:r>'\nSet-ParagraphRuns $paraA2 $runsA2\n"}
